$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, pushing existing rows 6-15 down to 7-16.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly entry.
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C6").Value = "Los Lagos"
$ws.Range("D6").Value = 44482
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 300000000
$ws.Range("G6").Value = "Espárragos"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 72
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 2000
$ws.Range("N6").Value = "$/kilo"
$ws.Range("O6").Value = "Provincia de Linares"
$ws.Range("P6").Value = 2000
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = "Hortaliza"
